$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - North America
$ws.Range("B2").Value = 127514226
$ws.Range("D2").Value = 1642146
$ws.Range("F2").Value = 123371297
$ws.Range("G2").Value = 257
$ws.Range("H2").Value = 2500783
$ws.Range("I2").Value = 6121

# Row 3 - Asia
$ws.Range("B3").Value = 220006270
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1548725
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = 203301226
$ws.Range("G3").Value = 64155
$ws.Range("H3").Value = 15156319
$ws.Range("I3").Value = 15135

# Row 4 - Europe
$ws.Range("B4").Value = 249785633
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 2069477
$ws.Range("F4").Value = 245928953
$ws.Range("G4").Value = 580
$ws.Range("H4").Value = 1787203
$ws.Range("I4").Value = 5443

# Row 5 - South America
$ws.Range("B5").Value = 68933479
$ws.Range("D5").Value = 1359505
$ws.Range("F5").Value = 66510589
$ws.Range("H5").Value = 1063385

# Row 6 - Australia/Oceania
$ws.Range("B6").Value = 14608586
$ws.Range("D6").Value = 30035
$ws.Range("F6").Value = 14469937
$ws.Range("H6").Value = 108614
$ws.Range("I6").Value = 49

# Row 7 - Africa
$ws.Range("B7").Value = 12836690
$ws.Range("D7").Value = 258825
$ws.Range("F7").Value = 12088094
$ws.Range("H7").Value = 489771
